$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the whole Price column as Text so that purely numeric-looking
# strings (e.g. "214.93") are stored as text rather than being auto-converted
# to numbers, matching the source data which keeps these as strings.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.956.79'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.676.44'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '214.93'
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").Value = '  -2.22%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = '20.95'
$ws.Range("E10").Value = '  +3.85%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '1.913.19'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '1.686.87'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("D16").Value = '65.73'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '8.19'
$ws.Range("E17").Value = '  +5.54%  '
$ws.Range("D18").Value = '26.969.97'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '235.92'
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D22").Value = '4.43'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("E24").Value = '  -4.41%  '
$ws.Range("D25").Value = '146.89'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").Value = '7.22'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").Value = '16.04'
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("E28").Value = '  -2.93%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("D33").Value = '1.488.32'
$ws.Range("E33").Value = '  +2.52%  '
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("E35").Value = '  +4.86%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  +3.71%  '
$ws.Range("D38").Value = '0.0174'
$ws.Range("E38").Value = '  +3.25%  '
$ws.Range("D39").Value = '0.910'
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("E40").Value = '  +4.74%  '
$ws.Range("D41").Value = '5.75'
$ws.Range("E41").Value = '  -5.03%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '67.38'
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("D45").Value = '1.818.45'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = '0.780'
$ws.Range("D47").Value = '90.46'
$ws.Range("D49").Value = '0.103'
$ws.Range("E49").Value = '  +2.64%  '

# Row 50 and 51 swap places (Cronos <-> EnergySwap)
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.76'
$ws.Range("E50").Value = '  +1.62%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0508'
$ws.Range("E51").Value = '  +0.28%  '

# Remove the temporary text formatting so the cell styles match the original
# workbook (which has no explicit style on these data cells).
$priceRange.ClearFormats()
